# Apply the ifoCAST sampling update to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 128: remove D128 (value previously -0.126881291408459) ---
$ws.Range("D128").Value = $null

# --- Rows 129-132: re-rounded D values ---
$ws.Range("D129").Value = 0.092077468
$ws.Range("D130").Value = 0.103503504
$ws.Range("D131").Value = 0.02315905499999998
$ws.Range("D132").Value = 0.08488154800000003

# --- Rows 133-139: re-rounded C values ---
$ws.Range("C133").Value = -0.09361257099999998
$ws.Range("C134").Value = -0.835520559
$ws.Range("C135").Value = -0.456828469
$ws.Range("C136").Value = -0.169753125
$ws.Range("C137").Value = -0.1421088
$ws.Range("C138").Value = -0.382958726
$ws.Range("C139").Value = -0.29709302

# --- Rows 136-139: new D values ---
$ws.Range("D136").Value = 1.068686474
$ws.Range("D137").Value = 1.064811887
$ws.Range("D138").Value = 0.874431902
$ws.Range("D139").Value = 0.8972184870000001

# --- Row 140: re-rounded B value and new C value ---
$ws.Range("B140").Value = -0.588418201
$ws.Range("C140").Value = 0.518578844

# --- New rows 141-145 ---
$ws.Range("A141").Value = "2025-07-25_diff"
$ws.Range("B141").Value = -0.6323316999999999
$ws.Range("C141").Value = 0.608342578

$ws.Range("A142").Value = "2025-08-07_diff"
$ws.Range("C142").Value = 0.430873215

$ws.Range("A143").Value = "2025-08-22_diff"
$ws.Range("C143").Value = 0.443590468

$ws.Range("A144").Value = "2025-08-25_diff"
$ws.Range("C144").Value = 0.319050665

$ws.Range("A145").Value = "2025-09-08_diff"
$ws.Range("C145").Value = 0.543389307

# Copy the existing column-A style (bold, bordered, centered) onto the new A cells
$ws.Range("A2").Copy()
$ws.Range("A141:A145").PasteSpecial(-4122)
$excel.CutCopyMode = $false
